# Microsite Education Script completed
# Adds new Help-Desk sprint-history rows to the AMSIN, BETA and AMS sheets,
# and fixes up the previously mis-formatted row 27 on the AMS sheet.

$wb = $excel.ActiveWorkbook

function Set-HistoryRow($ws, $row, $dateText, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    # Column A holds the run-date as literal text (not a real date value).
    # A leading apostrophe forces text interpretation instead of Excel's
    # automatic date parsing, and resetting the Style afterwards keeps the
    # effective format plain/General (matching the rest of the column).
    $ws.Cells.Item($row, 1).Value = "'" + $dateText
    $ws.Cells.Item($row, 1).Style = "Normal"

    # Column B holds the real run date/time serial, formatted the same way
    # as every other row in the table.
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $runTime

    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 4).Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 5).Style = "Normal"

    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 6).Style = "Normal"

    $ws.Cells.Item($row, 7).Value = $timeTaken
    $ws.Cells.Item($row, 7).Style = "Normal"
}

# ---------------------------------------------------------------------
# AMSIN sheet - append sprint 165 (3 cycles) and sprint 166 (2 cycles)
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Set-HistoryRow $wsAmsin 45 "2022-08-02" 44775.64484304399 "165_fstcycle"  124 122 2 2.1
Set-HistoryRow $wsAmsin 46 "2022-08-03" 44776.65965795139 "165_scndcycle" 124 121 3 2.05
Set-HistoryRow $wsAmsin 47 "2022-08-04" 44777.38325747685 "165_finalrun"  124 123 1 1.86
Set-HistoryRow $wsAmsin 48 "2022-08-22" 44795.65727164352 "166fstcycle"   124 122 2 2.07
Set-HistoryRow $wsAmsin 49 "2022-08-23" 44796.89554321759 "166cyclescnd"  124 122 2 1.82

# ---------------------------------------------------------------------
# BETA sheet - append the matching beta rows for sprint 165 and 166
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
Set-HistoryRow $wsBeta 25 "2022-08-04" 44777.5532071412  "165beta"  124 124 0 1.62
Set-HistoryRow $wsBeta 26 "2022-08-24" 44797.51717090278 "166_beta" 124 124 0 1.79

# ---------------------------------------------------------------------
# AMS sheet - fix up row 27 (was missing styling + had a stale run time),
# then append the "live" rows for sprint 165 and 166.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(27, 1).Style = "Normal"
$wsAms.Cells.Item(27, 2).Value = 44756.81026972222
$wsAms.Cells.Item(27, 3).Style = "Normal"
$wsAms.Cells.Item(27, 4).Style = "Normal"
$wsAms.Cells.Item(27, 5).Style = "Normal"
$wsAms.Cells.Item(27, 6).Style = "Normal"
$wsAms.Cells.Item(27, 7).Style = "Normal"

Set-HistoryRow $wsAms 28 "2022-08-04" 44777.8072408912 "165_live" 124 124 0 1.83

# Row 29 mirrors how row 27 originally looked before the fix above -
# values only, no explicit styling on the text/number cells, just the
# date/time format on column B. Resetting the Style after the text
# assignment clears the "quote prefix" flag the apostrophe trick leaves
# behind, so the effective format is plain/General like the rest of the
# un-styled cells in this table.
$wsAms.Cells.Item(29, 1).Value = "'2022-08-24"
$wsAms.Cells.Item(29, 1).Style = "Normal"
$wsAms.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Cells.Item(29, 2).Value = 44797.91147664575
$wsAms.Cells.Item(29, 3).Value = "166_live"
$wsAms.Cells.Item(29, 4).Value = 124
$wsAms.Cells.Item(29, 5).Value = 123
$wsAms.Cells.Item(29, 6).Value = 1
$wsAms.Cells.Item(29, 7).Value = 1.88
